$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-30 (A:BSL NO, B:BRAND, C:ISL NO, D:Item Name, E:UOM)
$rows = @(
  @(4, "Desodin", 1, "Desodin 60ml Syrup", "60 ml"),
  @(5, "Dinafex", 2, "Dinafex 180mg Tablet", "30's"),
  @(5, "Dinafex", 3, "Dinafex 120mg Tablet", "30's"),
  @(5, "Dinafex", 4, "Dinafex 60mg Tablet", "30's"),
  @(6, "Dorenta", 5, "Dorenta 50mg Tablet", "50's"),
  @(7, "Etorix", 6, "Etorix 60mg Tablet - 40's", "40's"),
  @(7, "Etorix", 7, "Etorix 120mg Tablet", "20's"),
  @(7, "Etorix", 8, "Etorix 90mg Tablet", "30's"),
  @(8, "Fenobac", 9, "Fenobac 100ml Syrup", "100ml"),
  @(9, "Flucloxin", 10, "Flucloxin 500mg Capsule", "30 's"),
  @(9, "Flucloxin", 11, "Flucloxin 500mg Capsule - 36's", "36 's"),
  @(10, "Geminox", 12, "Geminox 320mg Tablet - 8's", "8 's"),
  @(11, "Ketonic", 13, "Ketonic 30mg Injection", "5 's"),
  @(11, "Ketonic", 14, "Ketonic 30mg IM/IV Injection - 4's", "4's"),
  @(11, "Ketonic", 15, "Ketonic 10mg Tablet", "20's"),
  @(12, "Kynol", 16, "Kynol TR 200mg Capsule", "30 's"),
  @(12, "Kynol", 17, "Kynol D 25mg Tablet", "60 's"),
  @(12, "Kynol", 18, "Kynol TR 100mg Capsule", "50 's"),
  @(13, "Levomax", 19, "Levomax 750mg Tablet - 10's", "10 's"),
  @(17, "Naprox", 20, "Naprox Plus 500mg Tablet - 30's", "30 's"),
  @(19, "Oradin", 21, "Oradin Plus Tablet - 40's", "40 's"),
  @(20, "Osticare", 22, "Osticare Tablet 24's", "24's"),
  @(22, "Quinox", 23, "Quinox 100ml IV Infusion", "1's"),
  @(23, "Rupaday", 24, "Rupaday Oral Solution 60ml", "1's"),
  @(24, "Sk-Mox", 25, "Sk-Mox 500mg Capsule", "48 's"),
  @(35, "Zithrox", 26, "Zithrox 30ml Dry Suspension", "30ml"),
  @(35, "Zithrox", 27, "Zithrox 250mg Tablet - 6's", "6's"),
  @(35, "Zithrox", 28, "Zithrox 500mg Tablet", "6 's"),
  @(35, "Zithrox", 29, "Zithrox 15ml Suspension", "15 ml")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 2
  $row = $rows[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
}
